$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.713.41"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.600.00"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'211.57"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.248"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "1.824.73"
$ws.Range("D13").Value = "1.591.06"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'65.31"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "26.690.93"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'209.22"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'7.10"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'0.0523"
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "1.292.30"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'1.10"
$ws.Range("E39").Value = "  +20.03%  "
$ws.Range("D40").Value = "'0.824"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "'5.43"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'63.08"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").Value = "1.736.47"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "'0.0511"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +0.27%  "
